$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (intern_.in_parenthesis_replacement.R): fill in "x" across the columns
# that are already marked in neighbouring rows, except B15/D15/N15/O15 which
# stay empty because this item is blocked (see commit message).
$cols15 = @("C15","E15","F15","G15","H15","I15","J15","K15","L15","M15")
foreach ($addr in $cols15) {
    $ws.Range($addr).Value = "x"
}

# Row 19: remove the stray "x" mark in D19 (B19 stays as-is)
$ws.Range("D19").Clear()

# Update the remembered selection to match the last-edited cell
$ws.Range("M18").Select()
